# Insert a new data row at row 109 (pushing existing rows 109..210 down to
# 110..211) and populate it with the latest week's price data for
# Acelga / Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 109..210 down by one to make room for the new weekly record.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with this week's values. Columns that
# don't vary for this market/product combination (A, B, C, E, F, G, H, I, N,
# O, Q, R) keep the same values as the surrounding rows.
$ws.Range("A109").Value = 3
$ws.Range("B109").Value = "Femacal de La Calera"
$ws.Range("C109").Value = "Coquimbo"
$ws.Range("D109").Value = 44484
$ws.Range("E109").Value = 5
$ws.Range("F109").Value = 100112009
$ws.Range("G109").Value = "Acelga"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 300
$ws.Range("K109").Value = 2000
$ws.Range("L109").Value = 2200
$ws.Range("M109").Value = 2107
$ws.Range("N109").Value = "$/docena de atados (6 kilos)"
$ws.Range("O109").Value = "Provincia de Quillota"
$ws.Range("P109").Value = 351
$ws.Range("Q109").Value = 6
$ws.Range("R109").Value = "Hortaliza"
